# Update the fixed footer "Date Placeholder" text from 1/29/2025 to
# 1/30/2025 across the deck's slide master and every slide layout
# (stale date left over from Mesh testing).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq "1/29/2025") {
                $sh.TextFrame.TextRange.Text = "1/30/2025"
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}
